$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Copy the formatting used by the "status" column cells (B32:B38, style index 4 -
#    green fill + top alignment) down onto B40:B41, then stamp them "Done".
$ws.Range("B38").Copy() | Out-Null
$ws.Range("B40:B41").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("B40").Value = "Done"
$ws.Range("B41").Value = "Done"

# 2) A39: "Koersborden" -> "Koersborden/loknrs"
$ws.Range("A39").Value = "Koersborden/loknrs"

# 3) Update the view: scrolled one row further, and selection moved to A40.
$ws.Application.ActiveWindow.ScrollRow = 24
$ws.Range("A40").Select() | Out-Null
